$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two existing hyperlinks up front (the unitedWayHouston one in
# row 3 and the archived salvArmy one in row 2) before touching any cell
# content, so the Hyperlinks collection isn't mutated mid-edit.
$ws.Range("B3").Hyperlinks.Delete()
$ws.Range("B2").Hyperlinks.Delete()

# --- Row 2 (salvArmy): replace the archived Wayback Machine link with the
#     live Salvation Army URL, both as the displayed text and as the
#     hyperlink target. ---
$ws.Range("B2").Value = "https://centralusa.salvationarmy.org/midland/news/the-salvation-army-midland-division-opens-cooling-centers-with-missouri-and-illinois-to-battle-summer-heat/"
$ws.Hyperlinks.Add($ws.Range("B2"), "https://centralusa.salvationarmy.org/midland/news/the-salvation-army-midland-division-opens-cooling-centers-with-missouri-and-illinois-to-battle-summer-heat/")
$ws.Range("B2").Style = "Hyperlink"

# --- Row 3 (was unitedWayHouston): drop that site and replace it with the
#     NJ 211 cooling-centers entry that used to live in row 4. NJ is kept
#     as plain text (no hyperlink), but the cell keeps the leftover
#     Hyperlink formatting from the removed link. ---
$ws.Range("A3").Value = "NJ"
$ws.Range("B3").Value = "https://web.archive.org/web/20210701183644/https://www.nj211.org/nj-cooling-centers"

# --- Row 4: now empty (its data moved up into row 3). B4 keeps the
#     Hyperlink-style formatting even though it has no content. ---
$ws.Range("A4").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B4").Style = "Hyperlink"

# --- Selection moves to A4. ---
$ws.Range("A4").Select()
